# Liverpool_stats.xlsx update
#
# The sofascore export re-ordered two player rows (Dominik Szoboszlai now
# sorts ahead of Alexis Mac Allister) and refreshed the betting-markets
# derived stats, which on this sheet shows up as row 16 (Virgil van Dijk)
# picking up one extra totw appearance (totwAppearances, column I).
#
# Net effect on the worksheet data:
#   * Row 8 and Row 9 swap their entire contents (every data column) -
#     Dominik Szoboszlai's stat line moves from row 9 to row 8, and
#     Alexis Mac Allister's stat line moves from row 8 to row 9.
#   * Row 16, column I (totwAppearances) goes from 2 to 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A (1) through DL (116) hold this table's data. Column DL
# (goalsPrevented) is blank for both rows either way, so restrict the
# swap to A:DK (1:115) and leave that trailing empty cell untouched.
$firstCol = 1
$lastCol  = 115

$rowA = 8   # Alexis Mac Allister (before the edit)
$rowB = 9   # Dominik Szoboszlai (before the edit)

for ($col = $firstCol; $col -le $lastCol; $col++) {
    $cellA = $ws.Cells.Item($rowA, $col)
    $cellB = $ws.Cells.Item($rowB, $col)

    $valA = $cellA.Value2
    $valB = $cellB.Value2

    $cellA.Value2 = $valB
    $cellB.Value2 = $valA
}

# Virgil van Dijk's totwAppearances (column I) increases from 2 to 3.
$ws.Range("I16").Value2 = 3
